$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1168.6631
$ws.Range("I15").Value = 1168.6631
$ws.Range("K15").Value = 3505.9893
$ws.Range("M15").Value = -3336.9893
$ws.Range("H112").Value = 1270.697
$ws.Range("I112").Value = 1026
$ws.Range("J112").Value = 1314.3928
$ws.Range("K112").Value = 3078
$ws.Range("L112").Value = 3943.1784
$ws.Range("M112").Value = -1970
$ws.Range("N112").Value = -6159.178400000001
$ws.Range("H113").Value = 57465.5
$ws.Range("J113").Value = 1884.5714
$ws.Range("L113").Value = 1884.5714
$ws.Range("N113").Value = -8392.571400000001
$ws.Range("H129").Value = 2416.3286
$ws.Range("I129").Value = 5383.3
$ws.Range("J129").Value = 1229.54
$ws.Range("K129").Value = 16149.9
$ws.Range("L129").Value = 3688.62
$ws.Range("M129").Value = -11149.9
$ws.Range("N129").Value = -13688.62
$ws.Range("H137").Value = 1103.6034
$ws.Range("I137").Value = 1098.9333
$ws.Range("J137").Value = 1119.7693
$ws.Range("K137").Value = 3296.7999
$ws.Range("L137").Value = 3359.3079
$ws.Range("M137").Value = -746.7999
$ws.Range("N137").Value = -8459.3079
$ws.Range("H138").Value = 2220.5293
$ws.Range("I138").Value = 1626.6957
$ws.Range("J138").Value = 3462.182
$ws.Range("K138").Value = 4880.0871
$ws.Range("L138").Value = 10386.546
$ws.Range("M138").Value = 259.9129000000003
$ws.Range("N138").Value = -20666.546

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3249.34
$ws.Range("I32").Value = 2852.1475
$ws.Range("J32").Value = 10796
$ws.Range("K32").Value = 2852.1475
$ws.Range("L32").Value = 10796
$ws.Range("M32").Value = -2565.1475
$ws.Range("N32").Value = -11370
$ws.Range("H110").Value = 41754572
$ws.Range("I110").Value = 43569948
$ws.Range("K110").Value = 43569948
$ws.Range("M110").Value = -43567903
$ws.Range("H122").Value = 1131.7273
$ws.Range("I122").Value = 1947
$ws.Range("K122").Value = 5841
$ws.Range("M122").Value = -3391
$ws.Range("H132").Value = 3028.45
$ws.Range("I132").Value = 3155.4468
$ws.Range("K132").Value = 9466.340400000001
$ws.Range("M132").Value = -6936.340400000001
$ws.Range("H138").Value = 38311.6
$ws.Range("J138").Value = 38311.6
$ws.Range("L138").Value = 38311.6
$ws.Range("N138").Value = -48591.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 43725
$ws.Range("J59").Value = 43725
$ws.Range("L59").Value = 43725
$ws.Range("N59").Value = -45419
$ws.Range("H94").Value = 453.35
$ws.Range("I94").Value = 424.57895
$ws.Range("J94").Value = 1000
$ws.Range("K94").Value = 424.57895
$ws.Range("L94").Value = 1000
$ws.Range("M94").Value = 26.42104999999998
$ws.Range("N94").Value = -1902
$ws.Range("H107").Value = 41686130
$ws.Range("I107").Value = 66697284
$ws.Range("J107").Value = 866.6667
$ws.Range("K107").Value = 66697284
$ws.Range("L107").Value = 866.6667
$ws.Range("M107").Value = -66695364
$ws.Range("N107").Value = -4706.6667
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H134").Value = 2276.5874
$ws.Range("I134").Value = 1825.0176
$ws.Range("J134").Value = 6566.5
$ws.Range("K134").Value = 5475.052799999999
$ws.Range("L134").Value = 19699.5
$ws.Range("M134").Value = -2940.052799999999
$ws.Range("N134").Value = -24769.5
$ws.Range("H135").Value = 54500
$ws.Range("J135").Value = 54500
$ws.Range("L135").Value = 54500
$ws.Range("N135").Value = -64640
$ws.Range("H137").Value = 56921.668
$ws.Range("J137").Value = 56921.668
$ws.Range("L137").Value = 56921.668
$ws.Range("N137").Value = -67121.66800000001
$ws.Range("H140").Value = 66995
$ws.Range("J140").Value = 66995
$ws.Range("L140").Value = 66995
$ws.Range("N140").Value = -77355

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1123
$ws.Range("I16").Value = 757.25
$ws.Range("J16").Value = 2098.3333
$ws.Range("K16").Value = 757.25
$ws.Range("L16").Value = 2098.3333
$ws.Range("M16").Value = -470.25
$ws.Range("N16").Value = -2672.3333
$ws.Range("H31").Value = 40969.324
$ws.Range("I31").Value = 1775.6666
$ws.Range("J31").Value = 53567.285
$ws.Range("K31").Value = 1775.6666
$ws.Range("L31").Value = 53567.285
$ws.Range("M31").Value = -1480.6666
$ws.Range("N31").Value = -54157.285
$ws.Range("H34").Value = 40969.324
$ws.Range("I34").Value = 1775.6666
$ws.Range("J34").Value = 53567.285
$ws.Range("K34").Value = 1775.6666
$ws.Range("L34").Value = 53567.285
$ws.Range("M34").Value = -1573.6666
$ws.Range("N34").Value = -53971.285
$ws.Range("H86").Value = 3488.125
$ws.Range("I86").Value = 2666.6667
$ws.Range("J86").Value = 3981
$ws.Range("K86").Value = 2666.6667
$ws.Range("L86").Value = 3981
$ws.Range("M86").Value = -1543.6667
$ws.Range("N86").Value = -6227
$ws.Range("H89").Value = 3488.125
$ws.Range("I89").Value = 2666.6667
$ws.Range("J89").Value = 3981
$ws.Range("K89").Value = 13333.3335
$ws.Range("L89").Value = 19905
$ws.Range("M89").Value = -7717.333500000001
$ws.Range("N89").Value = -31137
$ws.Range("H94").Value = 1032.5333
$ws.Range("I94").Value = 800.8
$ws.Range("J94").Value = 1148.4
$ws.Range("K94").Value = 800.8
$ws.Range("L94").Value = 1148.4
$ws.Range("M94").Value = -349.8
$ws.Range("N94").Value = -2050.4
$ws.Range("H105").Value = 1400.0834
$ws.Range("I105").Value = 1498.3334
$ws.Range("J105").Value = 1301.8334
$ws.Range("K105").Value = 1498.3334
$ws.Range("L105").Value = 1301.8334
$ws.Range("M105").Value = 248.6666
$ws.Range("N105").Value = -4795.8334
$ws.Range("H107").Value = 3788.1562
$ws.Range("I107").Value = 9783.546
$ws.Range("J107").Value = 647.7143
$ws.Range("K107").Value = 9783.546
$ws.Range("L107").Value = 647.7143
$ws.Range("M107").Value = -7863.546
$ws.Range("N107").Value = -4487.7143
$ws.Range("H113").Value = 1123
$ws.Range("I113").Value = 757.25
$ws.Range("J113").Value = 2098.3333
$ws.Range("K113").Value = 757.25
$ws.Range("L113").Value = 2098.3333
$ws.Range("M113").Value = 1412.75
$ws.Range("N113").Value = -6438.3333
$ws.Range("H122").Value = 517.75
$ws.Range("I122").Value = 374.4
$ws.Range("J122").Value = 620.1429000000001
$ws.Range("K122").Value = 1123.2
$ws.Range("L122").Value = 1860.4287
$ws.Range("M122").Value = 1326.8
$ws.Range("N122").Value = -6760.4287

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 800
$ws.Range("I58").Value = 800
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 2400
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -2272
$ws.Range("N58").ClearContents()
$ws.Range("H107").Value = 267690.3
$ws.Range("J107").Value = 556364.5
$ws.Range("L107").Value = 1669093.5
$ws.Range("N107").Value = -1672933.5
$ws.Range("H122").Value = 435
$ws.Range("J122").Value = 372.5
$ws.Range("L122").Value = 3352.5
$ws.Range("N122").Value = -8252.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2142.9333
$ws.Range("I102").Value = 1545.4445
$ws.Range("K102").Value = 1545.4445
$ws.Range("M102").Value = 76.55549999999994
$ws.Range("H139").Value = 52369.43
$ws.Range("J139").Value = 52369.43
$ws.Range("L139").Value = 52369.43
$ws.Range("N139").Value = -62649.43

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2558.4119
$ws.Range("I122").Value = 2422.5386
$ws.Range("K122").Value = 7267.6158
$ws.Range("M122").Value = -4817.6158
$ws.Range("H132").Value = 2428
$ws.Range("I132").Value = 2394.762
$ws.Range("J132").Value = 2583.111
$ws.Range("K132").Value = 7184.286
$ws.Range("L132").Value = 7749.333
$ws.Range("M132").Value = -4654.286
$ws.Range("N132").Value = -12809.333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 286284.44
$ws.Range("I81").Value = 167332
$ws.Range("J81").Value = 999999
$ws.Range("K81").Value = 334664
$ws.Range("L81").Value = 1999998
$ws.Range("M81").Value = -333603
$ws.Range("N81").Value = -2002120
$ws.Range("H84").Value = 286284.44
$ws.Range("I84").Value = 167332
$ws.Range("J84").Value = 999999
$ws.Range("K84").Value = 1673320
$ws.Range("L84").Value = 9999990
$ws.Range("M84").Value = -1668016
$ws.Range("N84").Value = -10010598
$ws.Range("H100").Value = 1000000
$ws.Range("I100").Value = 1000000
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 2000000
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -1999459
$ws.Range("N100").ClearContents()
$ws.Range("H113").Value = 549.05554
$ws.Range("J113").Value = 584.0909
$ws.Range("L113").Value = 1752.2727
$ws.Range("N113").Value = -6092.2727
$ws.Range("H132").Value = 2473.8367
$ws.Range("I132").Value = 2676.9443
$ws.Range("J132").Value = 1911.3846
$ws.Range("K132").Value = 8030.8329
$ws.Range("L132").Value = 5734.1538
$ws.Range("M132").Value = -5500.8329
$ws.Range("N132").Value = -10794.1538
$ws.Range("H136").Value = 627.2619
$ws.Range("I136").Value = 395.92105
$ws.Range("J136").Value = 2825
$ws.Range("K136").Value = 1187.76315
$ws.Range("L136").Value = 8475
$ws.Range("M136").Value = 1362.23685
$ws.Range("N136").Value = -13575
